$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "2 Query SQL"
$ws.Range("B1").Value = "Diagramma ER"

$ws.Range("B2").Select()
